# Helper: set a cell's value while forcing it to be stored as TEXT
# (prevents Excel's automatic "looks like a number" -> numeric conversion),
# then strips the left-over NumberFormat by re-applying the format copied
# from a pristine, never-touched donor cell so no stray style index is
# left applied on the cell (matches the un-styled "t=inlineStr" cells in
# the original workbook).
function Set-TextValue {
    param($ws, $cellAddr, [string]$text, $donorAddr)
    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $text
    $ws.Range($donorAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q3" sheet, positioned right after "总计"
#    and right before the existing "2022-Q2" sheet (mirrors how a new
#    quarterly tab would be added by copying the prior quarter's sheet
#    as a template).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Pristine, never-reformatted scratch cell used purely as a "format
# donor" so PasteSpecial(xlPasteFormats) can restore the default
# (unstyled) look on every text cell we touch below.
$donor = "Z100"

# Fill the new "2022-Q3" sheet with its own data (overwriting the
# template values copied from "2022-Q2").
Set-TextValue $q3 "D2" "11.73" $donor
Set-TextValue $q3 "E2" "93.96" $donor
Set-TextValue $q3 "F2" "3.86"  $donor
Set-TextValue $q3 "G2" "0.4528" $donor
$q3.Range("H2").Value = 9

Set-TextValue $q3 "D3" "11.73" $donor
Set-TextValue $q3 "E3" "93.96" $donor
Set-TextValue $q3 "F3" "3.86"  $donor
Set-TextValue $q3 "G3" "0.4528" $donor
$q3.Range("H3").Value = 9

Set-TextValue $q3 "D4" "5.92" $donor
Set-TextValue $q3 "E4" "93.96" $donor
Set-TextValue $q3 "F4" "3.86"  $donor
Set-TextValue $q3 "G4" "0.2285" $donor
$q3.Range("H4").Value = 9

Set-TextValue $q3 "D5" "-11.74" $donor
Set-TextValue $q3 "E5" "93.96" $donor
Set-TextValue $q3 "F5" "3.86"  $donor
Set-TextValue $q3 "G5" "-0.4532" $donor
$q3.Range("H5").Value = 9

# Fund codes / names for the new quarter sheet.
Set-TextValue $q3 "B2" "006679" $donor
Set-TextValue $q3 "C2" "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 A" $donor

Set-TextValue $q3 "B3" "162719" $donor
Set-TextValue $q3 "C3" "广发道琼斯美国石油开发与生产指数（QDII-LOF）A" $donor

Set-TextValue $q3 "B4" "006680" $donor
Set-TextValue $q3 "C4" "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 C" $donor

Set-TextValue $q3 "B5" "004243" $donor
Set-TextValue $q3 "C5" "广发道琼斯美国石油开发与生产指数（QDII-LOF）C" $donor

# Clear the scratch donor cell so it doesn't leave a stray value behind.
$q3.Range($donor).Clear()

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: every existing row shifts down
#    one quarter, and a new row is appended for the oldest quarter that
#    drops out of the rolling window.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.68

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 2.37

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 1.38

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 2.06

$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.29

# New row 7, with the row-6 formatting (border/alignment/font) carried
# down onto column A just like it is for every other data row.
$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = "2021-Q2"
$total.Cells.Item(7, 3).Value = 4
$total.Cells.Item(7, 4).Value = 1.04

$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. The Sheet.Copy() call above made the brand-new "2022-Q3" tab the
#    active one. Restore "2021-Q2" (the sheet that was active/selected
#    in the original workbook) as the selected tab so that tab-selection
#    state is not dragged onto the newly inserted sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
